$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 833.3333  # H19
$ws.Cells.Item(19, 9).Value = 682.7273  # I19
$ws.Cells.Item(19, 10).Value = 999  # J19
$ws.Cells.Item(19, 11).Value = 682.7273  # K19
$ws.Cells.Item(19, 12).Value = 999  # L19
$ws.Cells.Item(19, 13).Value = -507.7273  # M19
$ws.Cells.Item(19, 14).Value = -1349  # N19
$ws.Cells.Item(70, 8).Value = 1621.4584  # H70
$ws.Cells.Item(70, 9).Value = 958  # I70
$ws.Cells.Item(70, 10).Value = 2095.3572  # J70
$ws.Cells.Item(70, 11).Value = 2874  # K70
$ws.Cells.Item(70, 12).Value = 6286.071599999999  # L70
$ws.Cells.Item(70, 13).Value = -2604  # M70
$ws.Cells.Item(70, 14).Value = -6826.071599999999  # N70
$ws.Cells.Item(73, 8).Value = 1621.4584  # H73
$ws.Cells.Item(73, 9).Value = 958  # I73
$ws.Cells.Item(73, 10).Value = 2095.3572  # J73
$ws.Cells.Item(73, 11).Value = 2874  # K73
$ws.Cells.Item(73, 12).Value = 6286.071599999999  # L73
$ws.Cells.Item(73, 13).Value = -1938  # M73
$ws.Cells.Item(73, 14).Value = -8158.071599999999  # N73
$ws.Cells.Item(93, 8).Value = 35000  # H93
$ws.Cells.Item(93, 10).Value = 35000  # J93
$ws.Cells.Item(93, 12).Value = 35000  # L93
$ws.Cells.Item(93, 14).Value = -39992  # N93
$ws.Cells.Item(137, 8).Value = 2565915.5  # H137
$ws.Cells.Item(137, 9).Value = 1978.16  # I137
$ws.Cells.Item(137, 10).Value = 7144375  # J137
$ws.Cells.Item(137, 11).Value = 5934.48  # K137
$ws.Cells.Item(137, 12).Value = 21433125  # L137
$ws.Cells.Item(137, 13).Value = -3384.48  # M137
$ws.Cells.Item(137, 14).Value = -21438225  # N137

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7097.58  # H32
$ws.Cells.Item(32, 9).Value = 7176.044  # I32
$ws.Cells.Item(32, 11).Value = 7176.044  # K32
$ws.Cells.Item(32, 13).Value = -6889.044  # M32
$ws.Cells.Item(61, 8).Value = 3336.6667  # H61
$ws.Cells.Item(61, 9).Value = 2390.5  # I61
$ws.Cells.Item(61, 10).Value = 3809.75  # J61
$ws.Cells.Item(61, 11).Value = 2390.5  # K61
$ws.Cells.Item(61, 12).Value = 3809.75  # L61
$ws.Cells.Item(61, 13).Value = -2178.5  # M61
$ws.Cells.Item(61, 14).Value = -4233.75  # N61
$ws.Cells.Item(122, 8).Value = 1814.9048  # H122
$ws.Cells.Item(122, 9).Value = 1718.6875  # I122
$ws.Cells.Item(122, 10).Value = 2122.8  # J122
$ws.Cells.Item(122, 11).Value = 5156.0625  # K122
$ws.Cells.Item(122, 12).Value = 6368.400000000001  # L122
$ws.Cells.Item(122, 13).Value = -2706.0625  # M122
$ws.Cells.Item(122, 14).Value = -11268.4  # N122
$ws.Cells.Item(132, 8).Value = 5398.303  # H132
$ws.Cells.Item(132, 9).Value = 2009.4286  # I132
$ws.Cells.Item(132, 10).Value = 7895.3687  # J132
$ws.Cells.Item(132, 11).Value = 6028.2858  # K132
$ws.Cells.Item(132, 12).Value = 23686.1061  # L132
$ws.Cells.Item(132, 13).Value = -3498.2858  # M132
$ws.Cells.Item(132, 14).Value = -28746.1061  # N132
$ws.Cells.Item(135, 8).Value = 1000000000  # H135
$ws.Cells.Item(135, 10).Value = 1000000000  # J135
$ws.Cells.Item(135, 12).Value = 1000000000  # L135
$ws.Cells.Item(135, 14).Value = -1000010140  # N135
$ws.Cells.Item(136, 8).Value = 3336.6667  # H136
$ws.Cells.Item(136, 9).Value = 2390.5  # I136
$ws.Cells.Item(136, 10).Value = 3809.75  # J136
$ws.Cells.Item(136, 11).Value = 7171.5  # K136
$ws.Cells.Item(136, 12).Value = 11429.25  # L136
$ws.Cells.Item(136, 13).Value = -4621.5  # M136
$ws.Cells.Item(136, 14).Value = -16529.25  # N136

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(109, 8).Value = 26996.666  # H109
$ws.Cells.Item(109, 10).Value = 26996.666  # J109
$ws.Cells.Item(109, 12).Value = 26996.666  # L109
$ws.Cells.Item(109, 14).Value = -29770.666  # N109

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6537753.5  # H31
$ws.Cells.Item(31, 9).Value = 1319.561  # I31
$ws.Cells.Item(31, 10).Value = 33337134  # J31
$ws.Cells.Item(31, 11).Value = 1319.561  # K31
$ws.Cells.Item(31, 12).Value = 33337134  # L31
$ws.Cells.Item(31, 13).Value = -1024.561  # M31
$ws.Cells.Item(31, 14).Value = -33337724  # N31
$ws.Cells.Item(34, 8).Value = 6537753.5  # H34
$ws.Cells.Item(34, 9).Value = 1319.561  # I34
$ws.Cells.Item(34, 10).Value = 33337134  # J34
$ws.Cells.Item(34, 11).Value = 1319.561  # K34
$ws.Cells.Item(34, 12).Value = 33337134  # L34
$ws.Cells.Item(34, 13).Value = -1117.561  # M34
$ws.Cells.Item(34, 14).Value = -33337538  # N34
$ws.Cells.Item(135, 8).Value = 36348  # H135
$ws.Cells.Item(135, 10).Value = 36348  # J135
$ws.Cells.Item(135, 12).Value = 36348  # L135
$ws.Cells.Item(135, 14).Value = -46488  # N135

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 704.5  # H5
$ws.Cells.Item(5, 9).Value = 455.95  # I5
$ws.Cells.Item(5, 11).Value = 1367.85  # K5
$ws.Cells.Item(5, 13).Value = -1255.85  # M5
$ws.Cells.Item(6, 8).Value = 701.1177  # H6
$ws.Cells.Item(6, 9).Value = 135.92308  # I6
$ws.Cells.Item(6, 11).Value = 407.76924  # K6
$ws.Cells.Item(6, 13).Value = -294.76924  # M6
$ws.Cells.Item(9, 8).Value = 2397.2222  # H9
$ws.Cells.Item(9, 9).Value = 0  # I9
$ws.Cells.Item(9, 10).Value = 2397.2222  # J9
$ws.Cells.Item(9, 11).Value = 0  # K9
$ws.Cells.Item(9, 12).Value = 7191.6666  # L9
$ws.Cells.Item(9, 13).ClearContents()  # M9
$ws.Cells.Item(9, 14).Value = -7639.6666  # N9
$ws.Cells.Item(10, 8).Value = 493.45456  # H10
$ws.Cells.Item(10, 9).Value = 380.8889  # I10
$ws.Cells.Item(10, 10).Value = 1000  # J10
$ws.Cells.Item(10, 11).Value = 1142.6667  # K10
$ws.Cells.Item(10, 12).Value = 3000  # L10
$ws.Cells.Item(10, 13).Value = -1003.6667  # M10
$ws.Cells.Item(10, 14).Value = -3278  # N10
$ws.Cells.Item(11, 8).Value = 389.5625  # H11
$ws.Cells.Item(11, 9).Value = 112  # I11
$ws.Cells.Item(11, 10).Value = 482.08334  # J11
$ws.Cells.Item(11, 11).Value = 336  # K11
$ws.Cells.Item(11, 12).Value = 1446.25002  # L11
$ws.Cells.Item(11, 13).Value = -196  # M11
$ws.Cells.Item(11, 14).Value = -1726.25002  # N11
$ws.Cells.Item(122, 8).Value = 2710.7778  # H122
$ws.Cells.Item(122, 10).Value = 3094.9678  # J122
$ws.Cells.Item(122, 12).Value = 27854.7102  # L122
$ws.Cells.Item(122, 14).Value = -32754.7102  # N122
$ws.Cells.Item(135, 8).Value = 704.5  # H135
$ws.Cells.Item(135, 9).Value = 455.95  # I135
$ws.Cells.Item(135, 11).Value = 4103.55  # K135
$ws.Cells.Item(135, 13).Value = -1568.55  # M135

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(34, 8).Value = 16000  # H34
$ws.Cells.Item(34, 10).Value = 16000  # J34
$ws.Cells.Item(34, 12).Value = 16000  # L34
$ws.Cells.Item(34, 14).Value = -16536  # N34
$ws.Cells.Item(74, 8).Value = 0  # H74
$ws.Cells.Item(74, 10).Value = 0  # J74
$ws.Cells.Item(74, 12).Value = 0  # L74
$ws.Cells.Item(74, 14).ClearContents()  # N74
$ws.Cells.Item(76, 8).Value = 16000  # H76
$ws.Cells.Item(76, 10).Value = 16000  # J76
$ws.Cells.Item(76, 12).Value = 16000  # L76
$ws.Cells.Item(76, 14).Value = -16630  # N76
$ws.Cells.Item(77, 8).Value = 0  # H77
$ws.Cells.Item(77, 10).Value = 0  # J77
$ws.Cells.Item(77, 12).Value = 0  # L77
$ws.Cells.Item(77, 14).ClearContents()  # N77
$ws.Cells.Item(79, 8).Value = 16000  # H79
$ws.Cells.Item(79, 10).Value = 16000  # J79
$ws.Cells.Item(79, 12).Value = 16000  # L79
$ws.Cells.Item(79, 14).Value = -18184  # N79
$ws.Cells.Item(80, 8).Value = 2794.1304  # H80
$ws.Cells.Item(80, 9).Value = 2692.0588  # I80
$ws.Cells.Item(80, 10).Value = 3083.3333  # J80
$ws.Cells.Item(80, 11).Value = 2692.0588  # K80
$ws.Cells.Item(80, 12).Value = 3083.3333  # L80
$ws.Cells.Item(80, 13).Value = -1694.0588  # M80
$ws.Cells.Item(80, 14).Value = -5079.3333  # N80
$ws.Cells.Item(83, 8).Value = 2794.1304  # H83
$ws.Cells.Item(83, 9).Value = 2692.0588  # I83
$ws.Cells.Item(83, 10).Value = 3083.3333  # J83
$ws.Cells.Item(83, 11).Value = 13460.294  # K83
$ws.Cells.Item(83, 12).Value = 15416.6665  # L83
$ws.Cells.Item(83, 13).Value = -8468.293999999998  # M83
$ws.Cells.Item(83, 14).Value = -25400.6665  # N83
$ws.Cells.Item(86, 8).Value = 0  # H86
$ws.Cells.Item(86, 10).Value = 0  # J86
$ws.Cells.Item(86, 12).Value = 0  # L86
$ws.Cells.Item(86, 14).ClearContents()  # N86
$ws.Cells.Item(89, 8).Value = 0  # H89
$ws.Cells.Item(89, 10).Value = 0  # J89
$ws.Cells.Item(89, 12).Value = 0  # L89
$ws.Cells.Item(89, 14).ClearContents()  # N89
$ws.Cells.Item(102, 8).Value = 2389.074  # H102
$ws.Cells.Item(102, 9).Value = 1412.4736  # I102
$ws.Cells.Item(102, 10).Value = 4708.5  # J102
$ws.Cells.Item(102, 11).Value = 1412.4736  # K102
$ws.Cells.Item(102, 12).Value = 4708.5  # L102
$ws.Cells.Item(102, 13).Value = 209.5264  # M102
$ws.Cells.Item(102, 14).Value = -7952.5  # N102

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(94, 8).Value = 19872  # H94
$ws.Cells.Item(94, 10).Value = 19872  # J94
$ws.Cells.Item(94, 12).Value = 19872  # L94
$ws.Cells.Item(94, 14).Value = -21224  # N94
$ws.Cells.Item(103, 8).Value = 0  # H103
$ws.Cells.Item(103, 10).Value = 0  # J103
$ws.Cells.Item(103, 12).Value = 0  # L103
$ws.Cells.Item(103, 14).ClearContents()  # N103
$ws.Cells.Item(115, 8).Value = 30000  # H115
$ws.Cells.Item(115, 10).Value = 30000  # J115
$ws.Cells.Item(115, 12).Value = 30000  # L115
$ws.Cells.Item(115, 14).Value = -32350  # N115

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 624.44446  # H107
$ws.Cells.Item(107, 9).Value = 624.44446  # I107
$ws.Cells.Item(107, 11).Value = 1873.33338  # K107
$ws.Cells.Item(107, 13).Value = 46.66661999999997  # M107
$ws.Cells.Item(132, 8).Value = 4907.0835  # H132
$ws.Cells.Item(132, 9).Value = 4981.1665  # I132
$ws.Cells.Item(132, 10).Value = 4833  # J132
$ws.Cells.Item(132, 11).Value = 14943.4995  # K132
$ws.Cells.Item(132, 12).Value = 14499  # L132
$ws.Cells.Item(132, 13).Value = -12413.4995  # M132
$ws.Cells.Item(132, 14).Value = -19559  # N132
